# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) timestamps on the
# zh-cn and de-de report sheets, per the regenerated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-28 11:04:11"
$wsZhCn.Range("G2").Value = "2016-01-28 11:04:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-28 11:04:24"
$wsDeDe.Range("G2").Value = "2016-01-28 11:05:20"
